# Add the "Week 4" meeting-minutes section to the end of the document.
# The document's final paragraph is an empty, centered, bold/red/32pt
# placeholder (the same style used for the "Week 2" / "Week 3" headings
# earlier in the file). We turn it into the "Week 4" heading and then
# append the usual per-member update paragraphs (Conor, Ahmed, Boyi,
# Ciaran, Oluwafirebami, Kevin) after it, mirroring the structure used
# for the previous weeks.

$d = $word.ActiveDocument

$lastParagraph = $d.Paragraphs.Last
$targetRange = $lastParagraph.Range

# Raw WordprocessingML for the new "Week 4" heading paragraph plus the
# eleven paragraphs that follow it (member headings, tabs, and the one
# filled-in update for Ahmed). Range.InsertXML replaces the contents of
# the addressed range with this markup, so this both turns the existing
# empty trailing paragraph into the "Week 4" heading and inserts all of
# the following paragraphs right before the section properties.
$weekFourXml = @'
<w:p><w:pPr><w:ind w:left="360"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Week 4</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Conor:</w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Ahmed:</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Finalised the screen and began work on a transition for the display of the program. Began work the dropdown box</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>Boyi:</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:tab/></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">Ciaran:  </w:t></w:r><w:r><w:tab/></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr></w:p><w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>Oluwafirebami</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">: </w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:ind w:left="360"/><w:rPr><w:b/><w:bCs/><w:color w:val="FF0000"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>Kevin:</w:t></w:r></w:p>
'@

$targetRange.InsertXML($weekFourXml)
